$d = $word.ActiveDocument

# Locate the "Requisitos" bullet-list paragraph (contains the LOM3215 line).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*LOM3215*" -and $t -like "*LOM3263*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

$vtab = [char]11
$pText = $target.Range.Text
$parts = $pText.Split($vtab)

# Compute absolute [start,end) offsets for each "line" (run text + its break char),
# skipping the trailing empty fragment after the last break (that's the paragraph mark).
$offsets = New-Object System.Collections.ArrayList
$cursor = $pStart
for ($i = 0; $i -lt $parts.Length - 1; $i++) {
    $lineLen = $parts[$i].Length + 1   # +1 for the w:br (vtab) char
    $lineStart = $cursor
    $lineEnd = $cursor + $lineLen
    [void]$offsets.Add(@{ Start = $lineStart; End = $lineEnd; Text = $parts[$i] })
    $cursor = $lineEnd
}

# Find the run whose text starts with "LOM3263"
$moveIdx = -1
for ($i = 0; $i -lt $offsets.Count; $i++) {
    if ($offsets[$i].Text -like "LOM3263*") {
        $moveIdx = $i
        break
    }
}

if ($moveIdx -ge 0) {
    $moveRange = $d.Range($offsets[$moveIdx].Start, $offsets[$moveIdx].End)
    $moveText = $moveRange.Text

    # Delete it from its current position.
    $moveRange.Delete()

    # Insert the moved text (run text + its line break) at the front of the paragraph.
    $insertAt = $d.Range($pStart, $pStart)
    $insertAt.InsertBefore($moveText)
}

Write-Output "done"
